$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.762082934379578
$ws.Range("B1").Value = 3.277971029281616
$ws.Range("C1").Value = 3.718083620071411
$ws.Range("D1").Value = 4.122349739074707
$ws.Range("E1").Value = 1.322736382484436
